# Updated symbol list (price / 1h-volume columns) to match the refreshed
# coinranking.com snapshot.
#
# D/E columns are stored as literal text in the workbook (e.g. "304.34",
# "-0.87%") rather than as numbers/percentages, so every new value is
# written with a leading apostrophe to force Excel to keep it as text
# instead of auto-converting it to a Number/Percentage cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'304.34"
$ws.Range("E2").Value = "'-0.87%"
# Row 3
$ws.Range("D3").Value = "'36.22"
$ws.Range("E3").Value = "'-2.25%"
# Row 4
$ws.Range("D4").Value = "'5.031"
$ws.Range("E4").Value = "'-0.22%"
# Row 5
$ws.Range("D5").Value = "'0.07857"
$ws.Range("E5").Value = "'-0.57%"
# Row 6
$ws.Range("E6").Value = "'-3.27%"
# Row 7
$ws.Range("D7").Value = "'7.959"
$ws.Range("E7").Value = "'-0.87%"
# Row 8
$ws.Range("D8").Value = "'0.9208"
$ws.Range("E8").Value = "'-0.94%"
# Row 9
$ws.Range("D9").Value = "'0.09525"
$ws.Range("E9").Value = "'-4.04%"
# Row 10
$ws.Range("D10").Value = "'0.1848"
$ws.Range("E10").Value = "'-2.15%"
# Row 11
$ws.Range("D11").Value = "'0.08797"
$ws.Range("E11").Value = "'1.18%"
# Row 12
$ws.Range("D12").Value = "'0.03606"
$ws.Range("E12").Value = "'0.08%"
# Row 13
$ws.Range("D13").Value = "'0.09910"
$ws.Range("E13").Value = "'-0.48%"
# Row 14
$ws.Range("D14").Value = "'0.001429"
$ws.Range("E14").Value = "'-3.92%"
# Row 15
$ws.Range("D15").Value = "'0.005696"
$ws.Range("E15").Value = "'0.18%"
# Row 16
$ws.Range("D16").Value = "'3.468"
$ws.Range("E16").Value = "'0.06%"
# Row 17
$ws.Range("D17").Value = "'4.141"
$ws.Range("E17").Value = "'2.65%"
# Row 18
$ws.Range("D18").Value = "'2.656"
$ws.Range("E18").Value = "'13.53%"
# Row 19
$ws.Range("E19").Value = "'-1.80%"
# Row 20
$ws.Range("E20").Value = "'2.22%"
# Row 21
$ws.Range("D21").Value = "'5.172"
$ws.Range("E21").Value = "'4.81%"
# Row 22
$ws.Range("E22").Value = "'2.40%"
# Row 23
$ws.Range("D23").Value = "'0.04567"
$ws.Range("E23").Value = "'-0.31%"
# Row 24
$ws.Range("D24").Value = "'0.001233"
$ws.Range("E24").Value = "'-1.41%"
# Row 25
$ws.Range("E25").Value = "'-8.90%"
# Row 26
$ws.Range("D26").Value = "'0.0001302"
$ws.Range("E26").Value = "'-6.97%"
# Row 27
$ws.Range("D27").Value = "'0.0004753"
$ws.Range("E27").Value = "'75.08%"
# Row 39
$ws.Range("D39").Value = "'0.01854"
$ws.Range("E39").Value = "'1.15%"
# Row 40
$ws.Range("D40").Value = "'0.04716"
$ws.Range("E40").Value = "'-1.44%"
# Row 41
$ws.Range("D41").Value = "'0.007783"
$ws.Range("E41").Value = "'-1.45%"
# Row 42
$ws.Range("D42").Value = "'0.1383"
$ws.Range("E42").Value = "'-2.51%"
# Row 43
$ws.Range("D43").Value = "'0.007728"
$ws.Range("E43").Value = "'2.59%"
# Row 44
$ws.Range("D44").Value = "'0.002211"
$ws.Range("E44").Value = "'1.02%"
# Row 45
$ws.Range("E45").Value = "'5.90%"
# Row 46
$ws.Range("D46").Value = "'0.00006367"
$ws.Range("E46").Value = "'1.18%"
# Row 47
$ws.Range("E47").Value = "'0.18%"
# Row 48
$ws.Range("E48").Value = "'0.29%"
# Row 49
$ws.Range("D49").Value = "'51.78"
$ws.Range("E49").Value = "'45.66%"
# Row 50
$ws.Range("D50").Value = "'0.001901"
$ws.Range("E50").Value = "'-29.22%"
# Row 51
$ws.Range("D51").Value = "'0.00002102"
$ws.Range("E51").Value = "'0.18%"
